$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 218, pushing the existing rows 218-236 down to 220-238.
$ws.Rows.Item(218).Insert()
$ws.Rows.Item(218).Insert()

# New row 218: Ají, Inferno, Primera (week of 2022-07-04)
$ws.Cells.Item(218,1).Value  = 8
$ws.Cells.Item(218,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(218,3).Value  = "Coquimbo"
$ws.Cells.Item(218,4).Value  = 44746
$ws.Cells.Item(218,5).Value  = 4
$ws.Cells.Item(218,6).Value  = 100112021
$ws.Cells.Item(218,7).Value  = "Ají"
$ws.Cells.Item(218,8).Value  = "Inferno"
$ws.Cells.Item(218,9).Value  = "Primera"
$ws.Cells.Item(218,10).Value = 500
$ws.Cells.Item(218,11).Value = 15000
$ws.Cells.Item(218,12).Value = 16000
$ws.Cells.Item(218,13).Value = 15500
$ws.Cells.Item(218,14).Value = "$/caja 12 kilos"
$ws.Cells.Item(218,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(218,16).Value = 1292
$ws.Cells.Item(218,17).Value = 12
$ws.Cells.Item(218,18).Value = "Hortaliza"

# New row 219: Ají, Inferno, Segunda (week of 2022-07-04)
$ws.Cells.Item(219,1).Value  = 8
$ws.Cells.Item(219,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(219,3).Value  = "Coquimbo"
$ws.Cells.Item(219,4).Value  = 44746
$ws.Cells.Item(219,5).Value  = 4
$ws.Cells.Item(219,6).Value  = 100112021
$ws.Cells.Item(219,7).Value  = "Ají"
$ws.Cells.Item(219,8).Value  = "Inferno"
$ws.Cells.Item(219,9).Value  = "Segunda"
$ws.Cells.Item(219,10).Value = 340
$ws.Cells.Item(219,11).Value = 10000
$ws.Cells.Item(219,12).Value = 11000
$ws.Cells.Item(219,13).Value = 10500
$ws.Cells.Item(219,14).Value = "$/caja 12 kilos"
$ws.Cells.Item(219,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(219,16).Value = 875
$ws.Cells.Item(219,17).Value = 12
$ws.Cells.Item(219,18).Value = "Hortaliza"
